# Applies the "WIP for scenarios and added questions file" edit to the
# "template scénarios" workbook:
#  - Fills in the header block (B1:B5) with the user-story metadata.
#  - Fills in the first three scenario rows (A8:C10) of the table.
#  - Widens columns A, B, C to fit the new text.
#  - Moves the active selection to C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block -----------------------------------------------------
$ws.Range("B1").Value = "PW-0001"
$ws.Range("B2").Value = "utilisateur/acheteur potentiel"
$ws.Range("B3").Value = "naviguer sur le site internet"
$ws.Range("B4").Value = "m'informer"
$ws.Range("B5").Value = "M "

# --- Scenario table rows ------------------------------------------------
$ws.Range("A8").Value = "je tappe l'url dans la barre de mon navigateur internet"
$ws.Range("C8").Value = "le serveur répond et envoye la page index.php a l'utilisateur"

$ws.Range("A9").Value = "je clique sur un tank qui m'intéresse"
$ws.Range("B9").Value = "une erreur quelconque se produit  côté site web/serveur"
$ws.Range("C9").Value = "le serveur affiche une page d'erreur comme quoi la ressource demandée n'as pas pu être chargée et envoyée"

$ws.Range("A10").Value = "je clique sur un autre tank qui m'intéresse "
$ws.Range("C10").Value = "le serveur répond, charge et envoye la ressource demand"

# --- Column widths -------------------------------------------------
# Target (authored) widths are 52.140625 / 51.140625 / 99.42578125 chars.
# ColumnWidth is quantized to the screen pixel grid on write, so feed it
# width-minus-padding (5/6 char) to land on the pixel closest to the
# authored value.
$ws.Columns.Item(1).ColumnWidth = 52.140625 - 5/6
$ws.Columns.Item(2).ColumnWidth = 51.140625 - 5/6
$ws.Columns.Item(3).ColumnWidth = 99.42578125 - 5/6

# --- Selection ----------------------------------------------------------
$ws.Range("C10").Select()
